# Practical 4 solutions - slide 15 title was missing the exercise number.
# "Exercise" -> "Exercise 2" (added as separate runs, matching how
# PowerPoint's TextRange.InsertAfter appends new runs rather than merging
# into the existing one).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$tr = $s.Shapes.Item(1).TextFrame.TextRange
$null = $tr.InsertAfter(" ")
$null = $tr.InsertAfter("2")
